$wb = $excel.ActiveWorkbook

# Add a new worksheet named "AFFIRM" and move it to the end of the tab
# order (after the existing "AFCEA" sheet).
$ws = $wb.Worksheets.Add()
$ws.Name = "AFFIRM"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $lastSheet)
# Re-fetch: Move() re-indexes the sheet collection, which can leave the
# previously-held $ws reference pointing at the wrong sheet.
$ws = $wb.Worksheets.Item("AFFIRM")

# Header row
$ws.Range("A1").Value = "Event Name"
$ws.Range("B1").Value = "Event Date"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Learn More"
$ws.Range("A1:D1").Font.Bold = $true

# Row 2 - Supply Chain Risk Management speaker series
$ws.Range("A2").Value = "Virtual Monthly Speaker Series: Supply Chain Risk Management in a Global World"
$ws.Range("B2").Value = "WedApril19"
$ws.Range("C2").Value = "Thought leaders from government and industry discuss the nature of complex software environments and methods to better manage and mitigate supply chain risks."
$ws.Range("D2").Value = "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/9p8B3mdL?sourceTypeId=Website&mode=Attendee"
$ws.Hyperlinks.Add($ws.Range("D2"), "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/9p8B3mdL?sourceTypeId=Website&mode=Attendee") | Out-Null

# Row 3 - Federal IT Career Workshop
$ws.Range("A3").Value = "2023 Federal IT Career Workshop"
$ws.Range("B3").Value = "TueApril25"
$ws.Range("C3").Value = "Join AFFIRM for this in-person workshop as you navigate your path through a career in government."
$ws.Range("D3").Value = "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/MrD67RZP?sourceTypeId=Website&mode=Attendee"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/MrD67RZP?sourceTypeId=Website&mode=Attendee") | Out-Null

# Row 4 - Metaverse speaker series
$ws.Range("A4").Value = "Virtual Monthly Speaker Series: Building Blocks of Government in the Metaverse"
$ws.Range("B4").Value = "WedMay17"
$ws.Range("C4").Value = "Our panel will discuss the pitfalls and the promise of the future of technology in the metaverse."
$ws.Range("D4").Value = "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/Zrnd7e0P?sourceTypeId=Website&mode=Attendee"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/Zrnd7e0P?sourceTypeId=Website&mode=Attendee") | Out-Null

# Row 5 - Emerging Technologies and AI speaker series
$ws.Range("A5").Value = "Virtual Monthly Speaker Series: Emerging Technologies and AI"
$ws.Range("B5").Value = "WedJune21"
$ws.Range("C5").Value = "Save the Date - event information coming soon!"
$ws.Range("D5").Value = "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/DpBbzVkp?sourceTypeId=Website&mode=Attendee"
$ws.Hyperlinks.Add($ws.Range("D5"), "https://associationforfederalinformationresourcesmanagementaffirm.growthzoneapp.com/ap/Events/Register/DpBbzVkp?sourceTypeId=Website&mode=Attendee") | Out-Null

# Hyperlinks.Add stamps its own (slightly different, but visually identical)
# cell style; re-apply the sheet's existing "Hyperlink" cell style so these
# cells reference the same style index the rest of the workbook uses.
$ws.Range("D2:D5").Style = "Hyperlink"

Write-Output "AFFIRM sheet added"
